$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.325.07'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.780.08'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.11%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9993'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9996'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5222'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +10.71%  '
$ws.Range('E8').Value = '  +5.60%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.49'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.86%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07394'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.78%  '
$ws.Range('E11').Value = '  +5.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9992'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.60'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.070'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.780.69'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.36%  '
$ws.Range('E16').Value = '  +2.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.55'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.37%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001047'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06419'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9994'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.73'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.859'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.402.64'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.32'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.076'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.48'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.13'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.50%  '
$ws.Range('E28').Value = '  +13.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.986.33'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '121.43'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.064'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.47%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09771'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.551'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.599'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02243'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05991'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '11.24'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.857'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.81%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6163'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2028'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.24%  '
$ws.Range('E41').Value = '  +2.44%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.117'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.152'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.98%  '
$ws.Range('E44').Value = '  +5.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5775'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.55%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.635'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '121.77'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.91%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.894'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.48%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.113'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.97%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06721'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.06'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.26%  '
